$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Weekly crime statistics table updates (rows 14-29) ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("G14").NumberFormat = "general"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("H14").NumberFormat = "general"
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 6
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -70
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 11
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 69
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = -20.689655172413
$ws.Range("L16").Value = 25.454545454545
$ws.Range("M16").Value = -19.767441860465
$ws.Range("N16").Value = -82.262210796915
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 72
$ws.Range("J17").Value = 99
$ws.Range("K17").Value = -27.272727272727
$ws.Range("L17").Value = 35.849056603773
$ws.Range("M17").Value = 33.333333333333
$ws.Range("N17").Value = -56.88622754491
$ws.Range("D18").Value = 13
$ws.Range("E18").Value = -69.230769230769
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -39.285714285714
$ws.Range("I18").Value = 122
$ws.Range("J18").Value = 158
$ws.Range("K18").Value = -22.784810126582
$ws.Range("L18").Value = -6.153846153846
$ws.Range("M18").Value = -24.223602484472
$ws.Range("N18").Value = -83.264746227709
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 350
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 10.416666666666
$ws.Range("I19").Value = 430
$ws.Range("J19").Value = 363
$ws.Range("K19").Value = 18.457300275482
$ws.Range("L19").Value = 69.291338582677
$ws.Range("M19").Value = 145.714285714286
$ws.Range("N19").Value = 91.111111111111
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 107
$ws.Range("J20").Value = 111
$ws.Range("K20").Value = -3.603603603603
$ws.Range("L20").Value = 40.78947368421
$ws.Range("M20").Value = 12.631578947368
$ws.Range("N20").Value = -81.709401709401
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 24
$ws.Range("F21").Value = 117
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = -12.686567164179
$ws.Range("I21").Value = 803
$ws.Range("J21").Value = 829
$ws.Range("K21").Value = -3.13630880579
$ws.Range("L21").Value = 39.895470383275
$ws.Range("M21").Value = 40.384615384615
$ws.Range("N21").Value = -61.88894162316
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 1
$ws.Range("F22").NumberFormat = '#,##0'
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("H22").Value = 0
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 8
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = -27.272727272727
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C23").NumberFormat = "general"
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J23").Value = 16
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 6.666666666666
$ws.Range("M23").Value = -20
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 70
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 4.477611940298
$ws.Range("I24").Value = 589
$ws.Range("J24").Value = 638
$ws.Range("K24").Value = -7.680250783699
$ws.Range("L24").Value = 15.717092337917
$ws.Range("M24").Value = 60.928961748633
$ws.Range("C25").Value = 5
$ws.Range("E25").Value = 25
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 155
$ws.Range("J25").Value = 178
$ws.Range("K25").Value = -12.921348314606
$ws.Range("L25").Value = 14.814814814814
$ws.Range("M25").Value = 9.929078014184
$ws.Range("D26").Value = 1
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -83.333333333333
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = -58.333333333333
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 31
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 63.157894736842
$ws.Range("L27").Value = 55
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J28").Value = 6
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J29").Value = 6
